$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$userAgent = "Mozilla/5.0 (Macintosh; Intel Mac OS X 10_15_7) AppleWebKit/537.36 (KHTML, like Gecko) Chrome/138.0.0.0 Safari/537.36"

# The stray empty cell previously dangling at I9 is dropped (row 9 now ends at H9).
$ws.Cells.Item(9, 9).ClearContents()

# New log rows appended starting at row 10
$rows = @(
    @("2025-12-11T04:52:35.466358", "sunil", 2, "GET", "/", 200, "127.0.0.1", $userAgent),
    @("2025-12-11T04:52:35.524407", "sunil", 2, "GET", "/favicon.ico", 404, "127.0.0.1", $userAgent),
    @("2025-12-11T04:52:38.886638", "sunil", 2, "GET", "/docs", 200, "127.0.0.1", $userAgent),
    @("2025-12-11T04:52:39.039368", "sunil", 2, "GET", "/openapi.json", 200, "127.0.0.1", $userAgent),
    @("2025-12-11T04:53:58.648374", "sunil", 2, "POST", "/auth/secure-login", 200, "127.0.0.1", $userAgent)
)

$r = 10
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Last row (14) carries the empty trailing I cell, matching the prior pattern from row 9
$ws.Cells.Item(14, 9).Borders.LineStyle = 0
